$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Aanwezig / Afwezig / problemen / doen / goed / beter" answer
# cells from each scrum block, leaving the formatting in place, so the
# sheet is ready to be filled in for the next scrum.
$ws.Range("C6:C12").ClearContents()
$ws.Range("C14:C20").ClearContents()
$ws.Range("C23:C24").ClearContents()
$ws.Range("C29").ClearContents()

# Reset the active selection back to the top of the sheet.
$ws.Range("C3").Select()
